$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: duplicate row 44 (current/pre-edit state) down into rows 45 and 46 ---
# Row 44's columns K..AY already hold the exact values/format that the new
# rows 45 and 46 need, and F..J on row 44 are still "TBD" at this point,
# which also matches what rows 45/46 need there. So copy the whole row first,
# then patch the few cells that actually differ (formulas need to be
# re-applied because PasteSpecial(xlPasteAll) bakes formulas to static values).
$ws.Range("A44:AY44").Copy()
$ws.Range("A45:AY45").PasteSpecial(-4104)
$ws.Range("A46:AY46").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# --- Step 2: restore the formula cells on rows 45 and 46 (Copy/Paste froze them to values) ---
foreach ($r in 45, 46) {
    $ws.Range("U$r").Formula = "=S$r + T$r"
    $ws.Range("AE$r").Formula = "=1508.06553301511 + 0.00210606006752809 * (AM$r*AN$r*AO$r) / 5 * U$r"
    $ws.Range("AQ$r").Formula = "=FLOOR.MATH((AJ$r - AM$r) / 2)"
    $ws.Range("AR$r").Formula = "=FLOOR.MATH((AK$r - AN$r) / 2)"
    $ws.Range("AS$r").Formula = "=FLOOR.MATH((AL$r - AO$r) / 2)"
}

# --- Step 3: new session identifiers / goal / expectation text for rows 45 & 46 ---
# Order matters: each brand-new distinct string gets appended to the shared
# string table the first time it's written, so write them in the same order
# the author's workbook ended up with (230911-0, eval-scores text, 230911-1,
# higher-patience text) to line up the new shared-string indices.
$ws.Range("A45").Value = "230911-0"
$ws.Range("D45").Value = "eval scores will plateau marginally higher. I think the problem is the sample size being just too small."
$ws.Range("A46").Value = "230911-1"
$ws.Range("C45").Value = "higher patience, lower val frequency, let it run for longer (e.g. 72 instead of the usual 24 hours) dataset03/04, maybe the eval score will go up eventually"

# E45 stays "TBD" conceptually, but row 44's E44 is NOT "TBD" (it has the
# "(TBD wait till finish) running, ..." text), so explicitly reset E45/E46.
$ws.Range("E45").Value = "TBD"
$ws.Range("E46").Value = "TBD"
# Likewise C44/D44 are not "TBD", so reset C46/D46 (row 45 already got its
# own distinct C45/D45 text above). B46 stays "train3dunet" from the copy.
$ws.Range("C46").Value = "TBD"
$ws.Range("D46").Value = "TBD"

# --- Step 4: row 44 now has its own results filled in, so fix up F44:J44 ---
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = "NA"
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = "NA"

# --- Step 5: move the active cell/selection to where the author left off ---
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("J45").Select()
